# Weekly update: insert a new "Acelga" price record for Femacal de La Calera
# (Coquimbo / Provincia de Quillota) ahead of the existing row 171, shifting
# every subsequent record down by one row (old row 171 -> 172, ..., old row
# 188 -> 189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 171:188 down to 172:189, leaving a blank row 171 to fill in.
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new weekly record.
$ws.Range("A171").Value = 3
$ws.Range("B171").Value = "Femacal de La Calera"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").Value = 44449
$ws.Range("E171").Value = 5
$ws.Range("F171").Value = 100112009
$ws.Range("G171").Value = "Acelga"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 280
$ws.Range("K171").Value = 2000
$ws.Range("L171").Value = 2300
$ws.Range("M171").Value = 2129
$ws.Range("N171").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O171").Value = "Provincia de Quillota"
$ws.Range("P171").Value = 355
$ws.Range("Q171").Value = 6
$ws.Range("R171").Value = "Hortaliza"
